$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 3 new rows (17,18,19) by copying formatting from row 16, then set values for rows 10-19
$ws.Range("A16:M16").Copy()
$ws.Range("A17:M19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 10
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "Gaussian-Quadrature"
$ws.Cells.Item(10,3).Value = 1.623327862283947
$ws.Cells.Item(10,4).Value = 0.4783232991450093
$ws.Cells.Item(10,5).Value = 0.9523107647407856
$ws.Cells.Item(10,6).Value = 1.623327862283947
$ws.Cells.Item(10,7).Value = 0.6828143876837972
$ws.Cells.Item(10,8).Value = 1.35554530157708
$ws.Cells.Item(10,9).Value = 1.02730437444241
$ws.Cells.Item(10,10).Value = 0.4783232991450093
$ws.Cells.Item(10,11).Value = 0.7153170319428974
$ws.Cells.Item(10,12).Value = 1.169322447113422
$ws.Cells.Item(10,13).Value = 1.019937664978838

# Row 11
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "Spiral-90deg-10rot-5space"
$ws.Cells.Item(11,3).Value = 0.8714425786307027
$ws.Cells.Item(11,4).Value = 1.775106686671545
$ws.Cells.Item(11,5).Value = 1.092091832980529
$ws.Cells.Item(11,6).Value = 0.8714425786307027
$ws.Cells.Item(11,7).Value = 0.6437831501117769
$ws.Cells.Item(11,8).Value = 2.471385678643533
$ws.Cells.Item(11,9).Value = 0.8430965096567369
$ws.Cells.Item(11,10).Value = 1.775106686671545
$ws.Cells.Item(11,11).Value = 1.433599259826037
$ws.Cells.Item(11,12).Value = 1.15252091922837
$ws.Cells.Item(11,13).Value = 1.282817739449137

# Row 12
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "Spiral-90deg-15rot-5space"
$ws.Cells.Item(12,3).Value = 0.873160659192439
$ws.Cells.Item(12,4).Value = 1.779061393565964
$ws.Cells.Item(12,5).Value = 1.090787333581507
$ws.Cells.Item(12,6).Value = 0.873160659192439
$ws.Cells.Item(12,7).Value = 0.6455285897639559
$ws.Cells.Item(12,8).Value = 2.461900241651241
$ws.Cells.Item(12,9).Value = 0.8423637239002412
$ws.Cells.Item(12,10).Value = 1.779061393565964
$ws.Cells.Item(12,11).Value = 1.434924363573736
$ws.Cells.Item(12,12).Value = 1.154042511383087
$ws.Cells.Item(12,13).Value = 1.282133656942558

# Row 13
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "Spiral-90deg-10rot-3space"
$ws.Cells.Item(13,3).Value = 0.8712482041714945
$ws.Cells.Item(13,4).Value = 1.775199653706745
$ws.Cells.Item(13,5).Value = 1.091977760806366
$ws.Cells.Item(13,6).Value = 0.8712482041714945
$ws.Cells.Item(13,7).Value = 0.644183634166565
$ws.Cells.Item(13,8).Value = 2.471849147176957
$ws.Cells.Item(13,9).Value = 0.842128545353979
$ws.Cells.Item(13,10).Value = 1.775199653706745
$ws.Cells.Item(13,11).Value = 1.433588707256555
$ws.Cells.Item(13,12).Value = 1.152418455714025
$ws.Cells.Item(13,13).Value = 1.282764490897018

# Row 14
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "NoRotation-tilt60deg"
$ws.Cells.Item(14,3).Value = 0.5750160000000021
$ws.Cells.Item(14,4).Value = 1.232020000000006
$ws.Cells.Item(14,5).Value = 1.611004000000001
$ws.Cells.Item(14,6).Value = 0.5750160000000021
$ws.Cells.Item(14,7).Value = 0.3837960000000005
$ws.Cells.Item(14,8).Value = 3.303915999999993
$ws.Cells.Item(14,9).Value = 1.087603999999996
$ws.Cells.Item(14,10).Value = 1.232020000000006
$ws.Cells.Item(14,11).Value = 1.421512000000004
$ws.Cells.Item(14,12).Value = 0.998264000000003
$ws.Cells.Item(14,13).Value = 1.365559333333333

# Row 15
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "Rotation-NoTilt"
$ws.Cells.Item(15,3).Value = 0.01
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = 2.302249999999999
$ws.Cells.Item(15,6).Value = 0.01
$ws.Cells.Item(15,7).Value = 0
$ws.Cells.Item(15,8).Value = 4.637637499999999
$ws.Cells.Item(15,9).Value = 1.407774999999999
$ws.Cells.Item(15,10).Value = 0
$ws.Cells.Item(15,11).Value = 1.151125
$ws.Cells.Item(15,12).Value = 0.5805624999999998
$ws.Cells.Item(15,13).Value = 1.39294375

# Row 16
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "Rotation-60detTilt"
$ws.Cells.Item(16,3).Value = 0.4282489876480028
$ws.Cells.Item(16,4).Value = 0.4042059431936034
$ws.Cells.Item(16,5).Value = 1.739115263385601
$ws.Cells.Item(16,6).Value = 0.4282489876480028
$ws.Cells.Item(16,7).Value = 0.4110000596992008
$ws.Cells.Item(16,8).Value = 3.075290076262394
$ws.Cells.Item(16,9).Value = 1.22794206208
$ws.Cells.Item(16,10).Value = 0.4042059431936034
$ws.Cells.Item(16,11).Value = 1.071660603289602
$ws.Cells.Item(16,12).Value = 0.7499547954688025
$ws.Cells.Item(16,13).Value = 1.214300398711467

# Row 17
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "HexGrid-90degTilt5degRes"
$ws.Cells.Item(17,3).Value = 0.9735639125594217
$ws.Cells.Item(17,4).Value = 0.9941311660881583
$ws.Cells.Item(17,5).Value = 0.9975720190588698
$ws.Cells.Item(17,6).Value = 0.9735639125594217
$ws.Cells.Item(17,7).Value = 0.9857552454456368
$ws.Cells.Item(17,8).Value = 1.001209575530364
$ws.Cells.Item(17,9).Value = 0.9903891589145849
$ws.Cells.Item(17,10).Value = 0.9941311660881583
$ws.Cells.Item(17,11).Value = 0.9958515925735141
$ws.Cells.Item(17,12).Value = 0.9847077525664679
$ws.Cells.Item(17,13).Value = 0.9904368462661726

# Row 18
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = "HexGrid-90degTilt22p5degRes"
$ws.Cells.Item(18,3).Value = 1.228623931902858
$ws.Cells.Item(18,4).Value = 1.096688252106518
$ws.Cells.Item(18,5).Value = 0.9072947566675926
$ws.Cells.Item(18,6).Value = 1.228623931902858
$ws.Cells.Item(18,7).Value = 0.9850364823079195
$ws.Cells.Item(18,8).Value = 0.95902990685564
$ws.Cells.Item(18,9).Value = 0.9264444570442626
$ws.Cells.Item(18,10).Value = 1.096688252106518
$ws.Cells.Item(18,11).Value = 1.001991504387056
$ws.Cells.Item(18,12).Value = 1.115307718144957
$ws.Cells.Item(18,13).Value = 1.017186297814132

# Row 19
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(19,3).Value = 0.9705512491822333
$ws.Cells.Item(19,4).Value = 1.185537839581145
$ws.Cells.Item(19,5).Value = 0.93834062833232
$ws.Cells.Item(19,6).Value = 0.9705512491822333
$ws.Cells.Item(19,7).Value = 1.130380158708225
$ws.Cells.Item(19,8).Value = 0.7463180020601147
$ws.Cells.Item(19,9).Value = 0.9293919765687166
$ws.Cells.Item(19,10).Value = 1.185537839581145
$ws.Cells.Item(19,11).Value = 1.061939233956733
$ws.Cells.Item(19,12).Value = 1.016245241569483
$ws.Cells.Item(19,13).Value = 0.9834199757387924
